# Daily attendance processing - 2025-12-06 19:22:05
# Normalizes the "Recorded By" column (G) so multi-contributor session
# rows list contributors in reverse-recorded order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $raw = $cell.Value2

    if ($raw -ne $null -and $raw -ne "") {
        $parts = $raw -split ", "
        if ($parts.Count -gt 1) {
            $reversed = $parts[($parts.Count - 1)..0]
            $newVal = [string]::Join(", ", $reversed)
            $cell.Value = $newVal
        }
    }
}
